# Horarios actualizados Linea 141 - 696
# Update sheet "LP1912": header info rows + reordered/appended schedule rows (44-215)
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Cells.Item(2, 1).Value = "Última actualización: 12:42:00"
$ws1.Cells.Item(3, 1).Value = "Total filas: 210"

$sheet1Rows = @(
    @("06:01:37","07:31","16_SANTA ANA",90,"LP1912"),
    @("06:01:37","07:31","11_ETCHEVERRY",90,"LP1912"),
    @("07:31:43","07:31","10_OLMOS",0,"LP1912"),
    @("06:44:15","07:32","11_ETCHEVERRY",48,"LP1912"),
    @("06:01:37","07:32","84_COLONIA URQUIZA-ESC 49",91,"LP1912"),
    @("06:01:37","07:36","27_EL RETIRO",95,"LP1912"),
    @("06:44:15","07:37","27_EL RETIRO",53,"LP1912"),
    @("06:01:37","07:39","10_OLMOS",98,"LP1912"),
    @("07:31:43","07:43","11_ETCHEVERRY",12,"LP1912"),
    @("06:59:44","07:45","11_ETCHEVERRY",46,"LP1912"),
    @("06:01:37","07:47","14_ABASTO",106,"LP1912"),
    @("06:44:15","07:48","14_ABASTO",64,"LP1912"),
    @("06:01:37","07:51","215D_EL PATO",110,"LP1912"),
    @("07:31:43","07:55","14_ABASTO",24,"LP1912"),
    @("07:31:43","07:57","16_SANTA ANA",26,"LP1912"),
    @("07:31:43","07:59","23_HERNANDEZ",28,"LP1912"),
    @("07:31:43","08:00","23_HERNANDEZ",29,"LP1912"),
    @("07:31:43","08:03","11_ETCHEVERRY",32,"LP1912"),
    @("07:57:27","08:04","16_SANTA ANA",7,"LP1912"),
    @("06:44:15","08:04","23_HERNANDEZ",80,"LP1912"),
    @("06:59:44","08:05","23_HERNANDEZ",66,"LP1912"),
    @("06:44:15","08:12","15_ABASTO",88,"LP1912"),
    @("07:57:27","08:13","10_OLMOS",16,"LP1912"),
    @("06:44:15","08:21","26_HERNANDEZ",97,"LP1912"),
    @("08:20:43","08:21","215B_EL PATO",1,"LP1912"),
    @("06:59:44","08:22","16_P MOR-SANTA ANA",83,"LP1912"),
    @("06:44:15","08:23","16_P MOR-SANTA ANA",99,"LP1912"),
    @("06:44:15","08:23","215B_EL PATO",99,"LP1912"),
    @("07:57:27","08:26","23_HERNANDEZ",29,"LP1912"),
    @("06:44:15","08:27","84_COLONIA URQUIZA-ESC 49",103,"LP1912"),
    @("07:57:27","08:33","10_OLMOS",36,"LP1912"),
    @("08:20:43","08:37","23_HERNANDEZ",17,"LP1912"),
    @("06:44:15","08:42","81_EL PELIGRO",118,"LP1912"),
    @("07:31:43","08:43","14_ABASTO",72,"LP1912"),
    @("08:20:43","08:49","23_HERNANDEZ",29,"LP1912"),
    @("08:42:31","08:51","14_ABASTO",9,"LP1912"),
    @("07:31:43","08:52","81_EL PELIGRO",81,"LP1912"),
    @("08:42:31","08:53","10_OLMOS",11,"LP1912"),
    @("06:59:44","08:54","17_ROMERO",115,"LP1912"),
    @("08:55:44","08:55","17_ROMERO",0,"LP1912"),
    @("07:31:43","09:01","215A_EL PATO",90,"LP1912"),
    @("07:57:27","09:03","11_ETCHEVERRY",66,"LP1912"),
    @("07:31:43","09:10","16_P MOR-SANTA ANA",99,"LP1912"),
    @("08:42:31","09:10","16_SANTA ANA",28,"LP1912"),
    @("08:42:31","09:11","23_HERNANDEZ",29,"LP1912"),
    @("08:55:44","09:13","16_SANTA ANA",18,"LP1912"),
    @("07:31:43","09:16","27_EL RETIRO",105,"LP1912"),
    @("07:57:27","09:17","27_EL RETIRO",80,"LP1912"),
    @("07:57:27","09:18","81_EL PELIGRO",81,"LP1912"),
    @("07:31:43","09:21","26_HERNANDEZ",110,"LP1912"),
    @("08:20:43","09:22","16_SANTA ANA",62,"LP1912"),
    @("07:31:43","09:22","17_ROMERO",111,"LP1912"),
    @("07:57:27","09:23","17_ROMERO",86,"LP1912"),
    @("07:31:43","09:23","11_ETCHEVERRY",112,"LP1912"),
    @("08:55:44","09:24","23_HERNANDEZ",29,"LP1912"),
    @("09:31:25","09:31","23_HERNANDEZ",0,"LP1912"),
    @("07:57:27","09:32","15_ABASTO",95,"LP1912"),
    @("07:57:27","09:33","10_OLMOS",96,"LP1912"),
    @("08:42:31","09:34","23_HERNANDEZ",52,"LP1912"),
    @("08:20:43","09:41","215C_EL PATO",81,"LP1912"),
    @("09:31:25","09:42","16_SANTA ANA",11,"LP1912"),
    @("07:57:27","09:42","215C_EL PATO",105,"LP1912"),
    @("07:57:27","09:43","14_ABASTO",106,"LP1912"),
    @("08:55:44","09:52","15_ABASTO",57,"LP1912"),
    @("08:55:44","09:53","10_OLMOS",58,"LP1912"),
    @("09:31:25","09:59","23_HERNANDEZ",28,"LP1912"),
    @("09:31:25","10:02","17_ROMERO",31,"LP1912"),
    @("09:31:25","10:03","11_ETCHEVERRY",32,"LP1912"),
    @("08:42:31","10:07","10_OLMOS",85,"LP1912"),
    @("08:20:43","10:08","10_OLMOS",108,"LP1912"),
    @("08:20:43","10:12","15_ABASTO",112,"LP1912"),
    @("09:31:25","10:13","10_OLMOS",42,"LP1912"),
    @("08:42:31","10:21","26_HERNANDEZ",99,"LP1912"),
    @("09:31:25","10:22","16_SANTA ANA",51,"LP1912"),
    @("09:31:25","10:23","11_ETCHEVERRY",52,"LP1912"),
    @("08:42:31","10:26","215A_EL PATO",104,"LP1912"),
    @("10:30:21","10:31","16_SANTA ANA",1,"LP1912"),
    @("10:30:21","10:34","23_HERNANDEZ",4,"LP1912"),
    @("08:55:44","10:41","17_ROMERO",106,"LP1912"),
    @("09:31:25","10:42","17_ROMERO",71,"LP1912"),
    @("08:55:44","10:43","14_ABASTO",108,"LP1912"),
    @("10:30:21","10:46","16_SANTA ANA",16,"LP1912"),
    @("10:30:21","10:52","15_ABASTO",22,"LP1912"),
    @("10:30:21","10:53","10_OLMOS",23,"LP1912"),
    @("10:30:21","10:57","27_EL RETIRO",27,"LP1912"),
    @("10:30:21","10:59","23_HERNANDEZ",29,"LP1912"),
    @("09:31:25","10:59","27_EL RETIRO",88,"LP1912"),
    @("09:31:25","11:02","215C_EL PATO",91,"LP1912"),
    @("11:03:46","11:03","215C_EL PATO",0,"LP1912"),
    @("10:30:21","11:03","11_ETCHEVERRY",33,"LP1912"),
    @("11:03:46","11:04","16_SANTA ANA",1,"LP1912"),
    @("11:03:46","11:04","23_HERNANDEZ",1,"LP1912"),
    @("10:30:21","11:06","16_P MOR-167 Y 521",36,"LP1912"),
    @("10:30:21","11:11","10_OLMOS",41,"LP1912"),
    @("11:03:46","11:12","15_ABASTO",9,"LP1912"),
    @("09:31:25","11:17","16_P MOR-167 Y 521",106,"LP1912"),
    @("09:31:25","11:19","86_EST CHICA-ESC AGRARIA",108,"LP1912"),
    @("09:31:25","11:21","26_HERNANDEZ",110,"LP1912"),
    @("09:31:25","11:26","16_P MOR-SANTA ANA",115,"LP1912"),
    @("09:31:25","11:27","225_C ROCA-H SUR",116,"LP1912"),
    @("10:30:21","11:32","81_EL PELIGRO",62,"LP1912"),
    @("11:03:46","11:34","23_HERNANDEZ",31,"LP1912"),
    @("10:30:21","11:35","11_ETCHEVERRY",65,"LP1912"),
    @("10:30:21","11:39","16_P MOR-SANTA ANA",69,"LP1912"),
    @("11:03:46","11:41","17_ROMERO",38,"LP1912"),
    @("10:30:21","11:42","17_ROMERO",72,"LP1912"),
    @("11:03:46","11:43","10_OLMOS",40,"LP1912"),
    @("11:44:50","11:44","16_SANTA ANA",0,"LP1912"),
    @("11:44:50","11:44","16_P MOR-SANTA ANA",0,"LP1912"),
    @("10:30:21","11:48","10_OLMOS",78,"LP1912"),
    @("10:30:21","11:51","215B_EL PATO",81,"LP1912"),
    @("11:03:46","11:52","15_ABASTO",49,"LP1912"),
    @("10:30:21","11:54","15_ABASTO",84,"LP1912"),
    @("11:44:50","11:58","225_GOMEZ",14,"LP1912"),
    @("10:30:21","11:59","225_GOMEZ",89,"LP1912"),
    @("10:30:21","12:02","84_COLONIA URQUIZA-ESC 49",92,"LP1912"),
    @("11:44:50","12:04","23_HERNANDEZ",20,"LP1912"),
    @("12:04:10","12:04","16_SANTA ANA",0,"LP1912"),
    @("12:04:10","12:05","23_HERNANDEZ",1,"LP1912"),
    @("12:04:10","12:05","16_SANTA ANA",1,"LP1912"),
    @("10:30:21","12:06","16_P MOR-SANTA ANA",96,"LP1912"),
    @("11:03:46","12:06","14_ABASTO",63,"LP1912"),
    @("12:04:10","12:07","10_OLMOS",3,"LP1912"),
    @("12:04:10","12:07","16_P MOR-SANTA ANA",3,"LP1912"),
    @("11:03:46","12:09","10_OLMOS",66,"LP1912"),
    @("10:30:21","12:14","17_ROMERO",104,"LP1912"),
    @("11:44:50","12:16","16_SANTA ANA",32,"LP1912"),
    @("10:30:21","12:17","14_ABASTO",107,"LP1912"),
    @("11:03:46","12:20","14_ABASTO",77,"LP1912"),
    @("10:30:21","12:20","215A_EL PATO",110,"LP1912"),
    @("12:04:10","12:21","215A_EL PATO",17,"LP1912"),
    @("12:04:10","12:21","14_ABASTO",17,"LP1912"),
    @("10:30:21","12:21","26_HERNANDEZ",111,"LP1912"),
    @("11:03:46","12:31","17_ROMERO",88,"LP1912"),
    @("11:44:50","12:34","23_HERNANDEZ",50,"LP1912"),
    @("12:04:10","12:35","23_HERNANDEZ",31,"LP1912"),
    @("12:04:10","12:35","11_ETCHEVERRY",31,"LP1912"),
    @("11:03:46","12:36","27_EL RETIRO",93,"LP1912"),
    @("12:04:10","12:37","27_EL RETIRO",33,"LP1912"),
    @("11:44:50","12:37","11_ETCHEVERRY",53,"LP1912"),
    @("11:03:46","12:38","17_179 Y 38",95,"LP1912"),
    @("11:44:50","12:41","10_OLMOS",57,"LP1912"),
    @("12:42:00","12:44","10_OLMOS",2,"LP1912"),
    @("12:42:00","12:47","16_SANTA ANA",5,"LP1912"),
    @("11:03:46","12:48","11_ETCHEVERRY",105,"LP1912"),
    @("12:04:10","12:49","11_ETCHEVERRY",45,"LP1912"),
    @("11:44:50","13:02","15_ABASTO",78,"LP1912"),
    @("12:42:00","13:03","14_ABASTO",21,"LP1912"),
    @("12:42:00","13:05","23_HERNANDEZ",23,"LP1912"),
    @("12:04:10","13:07","16_P MOR-SANTA ANA",63,"LP1912"),
    @("12:42:00","13:08","10_OLMOS",26,"LP1912"),
    @("11:44:50","13:10","10_OLMOS",86,"LP1912"),
    @("11:44:50","13:13","215D_EL PATO",89,"LP1912"),
    @("12:04:10","13:14","215D_EL PATO",70,"LP1912"),
    @("12:42:00","13:15","11_ETCHEVERRY",33,"LP1912"),
    @("11:44:50","13:19","10_OLMOS",95,"LP1912"),
    @("12:04:10","13:20","10_OLMOS",76,"LP1912"),
    @("11:44:50","13:21","26_HERNANDEZ",97,"LP1912"),
    @("11:44:50","13:26","15_ABASTO",102,"LP1912"),
    @("11:44:50","13:26","14_ABASTO",102,"LP1912"),
    @("12:04:10","13:27","14_ABASTO",83,"LP1912"),
    @("12:42:00","13:36","15_ABASTO",54,"LP1912"),
    @("12:42:00","13:41","23_HERNANDEZ",59,"LP1912"),
    @("12:04:10","13:46","17_ROMERO",102,"LP1912"),
    @("12:42:00","13:47","16_SANTA ANA",65,"LP1912"),
    @("12:04:10","13:51","215A_EL PATO",107,"LP1912"),
    @("12:04:10","13:56","225_GOMEZ",112,"LP1912"),
    @("12:04:10","13:57","16_P MOR-167 Y 521",113,"LP1912"),
    @("12:42:00","14:04","17_ROMERO",82,"LP1912"),
    @("12:42:00","14:17","27_EL RETIRO",95,"LP1912"),
    @("12:42:00","14:20","215C_EL PATO",98,"LP1912"),
    @("12:42:00","14:21","26_HERNANDEZ",99,"LP1912")
)

$startRow1 = 44
for ($i = 0; $i -lt $sheet1Rows.Count; $i++) {
    $r = $startRow1 + $i
    $row = $sheet1Rows[$i]
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
}

# Update sheet "LP1912-215": header info rows + new appended row 28
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2, 1).Value = "Última actualización: 12:42:00"
$ws2.Cells.Item(3, 1).Value = "Total filas: 23"

$ws2.Cells.Item(28, 1).Value = "12:42:00"
$ws2.Cells.Item(28, 2).Value = "14:20"
$ws2.Cells.Item(28, 3).Value = "215C_EL PATO"
$ws2.Cells.Item(28, 4).Value = 98
$ws2.Cells.Item(28, 5).Value = "LP1912"

# Update sheet "6203-6173": header info rows + new appended row 36
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2, 1).Value = "Última actualización: 12:42:00"
$ws3.Cells.Item(3, 1).Value = "Total filas: 31"

$ws3.Cells.Item(36, 1).Value = "12:42:00"
$ws3.Cells.Item(36, 2).Value = "14:09"
$ws3.Cells.Item(36, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(36, 4).Value = 87
$ws3.Cells.Item(36, 5).Value = "L6173"
